{"js": "// Replace the 25 \"three-digit number divided by one-digit number\" problems\n// with their new values, one-to-one, by exact text match.\nconst replacements = [\n  [\"110\u00f72=\", \"315\u00f75=\"],\n  [\"731\u00f74=\", \"342\u00f77=\"],\n  [\"219\u00f79=\", \"860\u00f78=\"],\n  [\"759\u00f77=\", \"622\u00f75=\"],\n  [\"460\u00f77=\", \"865\u00f72=\"],\n  [\"461\u00f78=\", \"902\u00f79=\"],\n  [\"863\u00f74=\", \"302\u00f78=\"],\n  [\"613\u00f78=\", \"614\u00f79=\"],\n  [\"575\u00f77=\", \"131\u00f75=\"],\n  [\"326\u00f76=\", \"739\u00f72=\"],\n  [\"748\u00f77=\", \"869\u00f77=\"],\n  [\"962\u00f76=\", \"723\u00f74=\"],\n  [\"646\u00f75=\", \"557\u00f76=\"],\n  [\"673\u00f79=\", \"325\u00f79=\"],\n  [\"293\u00f79=\", \"536\u00f75=\"],\n  [\"641\u00f79=\", \"811\u00f73=\"],\n  [\"872\u00f74=\", \"326\u00f75=\"],\n  [\"154\u00f79=\", \"376\u00f77=\"],\n  [\"340\u00f72=\", \"918\u00f75=\"],\n  [\"685\u00f73=\", \"107\u00f76=\"],\n  [\"946\u00f75=\", \"698\u00f73=\"],\n  [\"960\u00f75=\", \"231\u00f79=\"],\n  [\"951\u00f73=\", \"443\u00f74=\"],\n  [\"701\u00f76=\", \"805\u00f78=\"],\n  [\"397\u00f79=\", \"223\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the 25 \"three-digit number divided by one-digit number\" problems\n# with their new values, one-to-one, by exact text match.\n$d = $word.ActiveDocument\n\n$pairs = [ordered]@{\n    \"110\u00f72=\" = \"315\u00f75=\"\n    \"731\u00f74=\" = \"342\u00f77=\"\n    \"219\u00f79=\" = \"860\u00f78=\"\n    \"759\u00f77=\" = \"622\u00f75=\"\n    \"460\u00f77=\" = \"865\u00f72=\"\n    \"461\u00f78=\" = \"902\u00f79=\"\n    \"863\u00f74=\" = \"302\u00f78=\"\n    \"613\u00f78=\" = \"614\u00f79=\"\n    \"575\u00f77=\" = \"131\u00f75=\"\n    \"326\u00f76=\" = \"739\u00f72=\"\n    \"748\u00f77=\" = \"869\u00f77=\"\n    \"962\u00f76=\" = \"723\u00f74=\"\n    \"646\u00f75=\" = \"557\u00f76=\"\n    \"673\u00f79=\" = \"325\u00f79=\"\n    \"293\u00f79=\" = \"536\u00f75=\"\n    \"641\u00f79=\" = \"811\u00f73=\"\n    \"872\u00f74=\" = \"326\u00f75=\"\n    \"154\u00f79=\" = \"376\u00f77=\"\n    \"340\u00f72=\" = \"918\u00f75=\"\n    \"685\u00f73=\" = \"107\u00f76=\"\n    \"946\u00f75=\" = \"698\u00f73=\"\n    \"960\u00f75=\" = \"231\u00f79=\"\n    \"951\u00f73=\" = \"443\u00f74=\"\n    \"701\u00f76=\" = \"805\u00f78=\"\n    \"397\u00f79=\" = \"223\u00f74=\"\n}\n\nforeach ($old in $pairs.Keys) {\n    $new = $pairs[$old]\n    $find = $d.Content.Find\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
